$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Лист1")

$ws.Rows("5:6").Delete()

$ws.Range("A25").Value = 5
$ws.Range("B25").Value = "Окружение"
$ws.Range("C25").Value = "Отключение электричества"
$ws.Range("D25").Value = "Сервер отключится, работа с клиентами временно прекратится"
$ws.Range("E25").Value = "Установка запасного электрогенератора (предполагается, что он был до внедрения проекта)"
$ws.Range("F25").Value = "Включить запасной генератор, если его нет официанты работают как до внедрения проекта"
$ws.Range("G25").Value = "Отключение электричества"
$ws.Range("H25").Value = 2
$ws.Range("I25").Value = 1
$ws.Range("J25").Formula = "=`$H25*`$I25"

$ws.Range("A26").Value = 6
$ws.Range("B26").Value = "Технические"
$ws.Range("C26").Value = "Поломка оборудования"
$ws.Range("D26").Value = "Прекращения работы ресторана"
$ws.Range("E26").Value = "Ежемесячный тех-осмотр"
$ws.Range("F26").Value = "Остановка работы, проведение технических работ"
$ws.Range("G26").Value = "Поломка оборудования"
$ws.Range("H26").Value = 2
$ws.Range("I26").Value = 1
$ws.Range("J26").Formula = "=`$H26*`$I26"

$ws.Range("A25:J26").HorizontalAlignment = -4108
$ws.Range("A25:J26").VerticalAlignment = -4108
$ws.Range("B25:J26").WrapText = $true

$ws.Rows("25").RowHeight = 75
$ws.Rows("26").RowHeight = 45

$ws.Range("C4").Select()
